$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($col, $row, $value) {
    $cell = $ws.Range("$col$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - BNB
Set-TextCell "D" 2 "245.42"
Set-TextCell "G" 2 "5"

# Row 3 - OKB
Set-TextCell "D" 3 "25.36"
Set-TextCell "G" 3 "5"

# Row 4 - HuobiToken
Set-TextCell "D" 4 "5.103"
Set-TextCell "G" 4 "5"

# Row 5 - Cronos
Set-TextCell "D" 5 "0.05570"
Set-TextCell "G" 5 "5"

# Row 6 - KuCoinToken
Set-TextCell "G" 6 "5"

# Row 7 - GateToken
Set-TextCell "D" 7 "3.018"
Set-TextCell "G" 7 "5"

# Row 8 - MXToken
Set-TextCell "D" 8 "0.8189"
Set-TextCell "G" 8 "5"

# Row 9 - FTXToken
Set-TextCell "D" 9 "0.8453"
Set-TextCell "G" 9 "5"

# Row 10 - was WazirX, now One
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D" 10 "0.0005951"
$ws.Range("E10").Value = "9OneONE"
Set-TextCell "G" 10 "5"

# Row 11 - was MandalaExchangeToken, now WazirX
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D" 11 "0.1340"
$ws.Range("E11").Value = "10WazirXWRX"
Set-TextCell "G" 11 "5"

# Row 12 - was LiechtensteinCryptoassetsExchange, now MandalaExchangeToken
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D" 12 "0.06949"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
Set-TextCell "G" 12 "5"

# Row 13 - BitrueCoin
Set-TextCell "D" 13 "0.02875"
Set-TextCell "G" 13 "5"

# Row 14 - BitMartToken
Set-TextCell "D" 14 "0.09377"
Set-TextCell "G" 14 "5"

# Row 15 - BitForexToken
Set-TextCell "G" 15 "5"

# Row 16 - was One, now TigerCash
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D" 16 "0.006149"
$ws.Range("E16").Value = "15TigerCashTCH"
Set-TextCell "G" 16 "5"

# Row 17 - was TigerCash, now LEO
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D" 17 "3.498"
$ws.Range("E17").Value = "16LEOLEO"
Set-TextCell "G" 17 "5"

# Row 18 - was LEO, now BTSEToken
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D" 18 "2.063"
$ws.Range("E18").Value = "17BTSETokenBTSE"
Set-TextCell "G" 18 "5"

# Row 19 - was BTSEToken, now BitpandaEcosystemToken
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell "D" 19 "0.3179"
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"
Set-TextCell "G" 19 "5"

# Row 20 - was BitpandaEcosystemToken, now LiechtensteinCryptoassetsExchange
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D" 20 "0.03191"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"
Set-TextCell "G" 20 "5"

# Row 21 - ProBitToken
Set-TextCell "G" 21 "5"

# Row 22 - MCDex
Set-TextCell "D" 22 "3.756"
Set-TextCell "G" 22 "5"

# Row 23 - CoinExToken
Set-TextCell "D" 23 "0.04733"
Set-TextCell "G" 23 "5"

# Row 24 - ZBToken
Set-TextCell "G" 24 "5"

# Row 25 - BitKan
Set-TextCell "G" 25 "5"

# Row 26 - HotbitToken
Set-TextCell "D" 26 "0.004637"
Set-TextCell "G" 26 "5"

# Row 27 - NitroEx
Set-TextCell "D" 27 "0.00009701"
$ws.Range("E27").Value = "26NitroExNTX"
Set-TextCell "G" 27 "5"

# Row 28 - UpBots
Set-TextCell "D" 28 "0.0001388"
Set-TextCell "G" 28 "5"

# Row 29 - Spectre.aiUtilityToken
Set-TextCell "G" 29 "5"

# Row 30 - LegolasExchange
Set-TextCell "G" 30 "5"

# Row 31 - BitZToken
Set-TextCell "G" 31 "5"

# Row 32 - Birake
Set-TextCell "G" 32 "5"

# Row 33 - NashExchange
Set-TextCell "G" 33 "5"

# Row 34 - AAXToken
Set-TextCell "G" 34 "5"

# Row 35 - CenX
Set-TextCell "G" 35 "5"

# Row 36 - BNIXToken
Set-TextCell "G" 36 "5"

# Row 37 - Polkally
Set-TextCell "G" 37 "5"

# Row 38 - Charli3
Set-TextCell "G" 38 "5"

# Row 39 - BlubitexToken
Set-TextCell "G" 39 "5"

# Row 40 - IDEX
Set-TextCell "D" 40 "0.03662"
Set-TextCell "G" 40 "5"

# Row 41 - was KickToken, now BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D" 41 "0.1350"
$ws.Range("E41").Value = "40BKEXTokenBKK"
Set-TextCell "G" 41 "5"

# Row 42 - was BKEXToken, now KickToken
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D" 42 "0.006191"
$ws.Range("E42").Value = "41KickTokenKICKBestin24h"
Set-TextCell "G" 42 "5"

# Row 43 - CEJI
Set-TextCell "D" 43 "0.002500"
Set-TextCell "G" 43 "5"

# Row 44 - LocalTraders
Set-TextCell "D" 44 "0.008339"
Set-TextCell "G" 44 "5"

# Row 45 - CoinLion
Set-TextCell "D" 45 "0.00005293"
Set-TextCell "G" 45 "5"

# Row 46 - Kangarootoken
Set-TextCell "G" 46 "5"

# Row 47 - CoinbaseStockToken
Set-TextCell "D" 47 "0.1500"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
Set-TextCell "G" 47 "5"

# Row 48 - BOLO
Set-TextCell "G" 48 "5"

# Row 49 - CryptobidCoin
Set-TextCell "D" 49 "0.00002100"
Set-TextCell "G" 49 "5"

# Row 50 - SpecialPowerGold
Set-TextCell "D" 50 "0.0002000"
Set-TextCell "G" 50 "5"

# Row 51 - DigiFinexToken
Set-TextCell "G" 51 "5"
